# "create customer F1 button added instead of webelement click"
# The createCustomer test now submits the form with the F1 key instead of
# clicking a web element, so the leftover/duplicate test data row on
# "Sheet3" is no longer needed and a bad bankIfscCode value gets fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Sheet3"

# Fix the typo'd IFSC code on the remaining "gow manickam" row (row 2).
$ws.Range("N2").Value = "UE121H1019"

# Remove the stale hyperlinks before we shuffle rows around - row deletion
# does not keep the Hyperlinks collection in sync with the shifted cells.
$ws.Hyperlinks.Delete()

# Row 3 (blank customerName test row) is obsolete - delete it; row 4
# ("reg cust4") shifts up to become the new row 3.
$ws.Rows.Item(3).Delete()

# Re-create the hyperlinks for the two data rows that remain.
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Pass@1234", "", "", "Pass@1234")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Pass@1234", "", "", "Pass@1234")

# This sheet becomes the active tab, with D2 selected.
$ws.Activate()
$ws.Range("D2").Select()
